$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "74.861.76"
$ws.Range("E2").Value = "  +1.66%  "
$ws.Range("D3").Value = "2.821.85"
$ws.Range("E3").Value = "  +8.25%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "188.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +3.33%  "
$ws.Range("E9").Value = "  -4.60%  "
$ws.Range("D10").Value = "2.818.78"
$ws.Range("E10").Value = "  +8.22%  "
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("E12").Value = "  +2.13%  "
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("D14").Value = "3.340.05"
$ws.Range("E14").Value = "  +7.31%  "
$ws.Range("D15").Value = "74.818.25"
$ws.Range("E15").Value = "  +1.68%  "
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.95"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.32%  "
$ws.Range("D18").Value = "2.808.16"
$ws.Range("E18").Value = "  +6.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "373.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.55%  "
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.99%  "
$ws.Range("D27").Value = "2.964.78"
$ws.Range("E27").Value = "  +8.47%  "
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.15%  "
$ws.Range("E30").Value = "  +8.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "511.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.22%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("E34").Value = "  -1.33%  "
$ws.Range("E35").Value = "  +2.69%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "20.06"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.119"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.41%  "
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "181.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +15.73%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.93%  "
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("E46").Value = "  +2.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.00%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0869"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.01%  "
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.58%  "
$ws.Range("E50").Value = "  +7.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.72%  "
